$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Columns E (5) and F (6), rows 2 through 7, hold percentage-like ratios
# that were stored as fractions (e.g. 0.88...) but should now be stored
# as the equivalent "already multiplied by 100" numbers (e.g. 88.31...).
for ($row = 2; $row -le 7; $row++) {
    foreach ($col in 5, 6) {
        $cell = $ws.Cells.Item($row, $col)
        $cell.Value2 = $cell.Value2 * 100
    }
}
